$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 178, pushing existing rows 178-205 down to 179-206.
$ws.Rows.Item(178).Insert()

# Populate the newly inserted row 178 with the new record.
$ws.Cells.Item(178, 1).Value = 3
$ws.Cells.Item(178, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(178, 3).Value = "Coquimbo"
$ws.Cells.Item(178, 4).Value = 44491
$ws.Cells.Item(178, 5).Value = 5
$ws.Cells.Item(178, 6).Value = 100112012
$ws.Cells.Item(178, 7).Value = "Espinaca"
$ws.Cells.Item(178, 8).Value = "Sin especificar"
$ws.Cells.Item(178, 9).Value = "Primera"
$ws.Cells.Item(178, 10).Value = 120
$ws.Cells.Item(178, 11).Value = 2500
$ws.Cells.Item(178, 12).Value = 2500
$ws.Cells.Item(178, 13).Value = 2500
$ws.Cells.Item(178, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(178, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(178, 16).Value = 833
$ws.Cells.Item(178, 17).Value = 3
$ws.Cells.Item(178, 18).Value = "Hortaliza"
